$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix mistake: IndID (column A) and Register (column B) were swapped for the
# "nra" register rows (rows 39-51). Swap A and B values back to the correct order.
for ($r = 39; $r -le 51; $r++) {
    $aCell = $ws.Cells.Item($r, 1)
    $bCell = $ws.Cells.Item($r, 2)
    $aVal = $aCell.Value2
    $bVal = $bCell.Value2
    $aCell.Value2 = $bVal
    $bCell.Value2 = $aVal
}

# Leave the cursor where the author ended up after making the fix.
$ws.Range("B63").Select() | Out-Null
